$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 12:42"

# Update country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B6").Value = 68594
$ws.Range("C6").Value = 383
$ws.Range("E6").Value = 67130
$ws.Range("B8").Value = 39502
$ws.Range("C8").Value = 2179
$ws.Range("E8").Value = 35733
$ws.Range("F9").Value = 2746
$ws.Range("B11").Value = 11435
$ws.Range("C11").Value = 538
$ws.Range("E11").Value = 11139
$ws.Range("B18").Value = 3212
$ws.Range("C18").Value = 128
$ws.Range("E18").Value = 3192
$ws.Range("B37").Value = 926
$ws.Range("C37").Value = 46
$ws.Range("E37").Value = 911
$ws.Range("B57").Value = 457
$ws.Range("C57").Value = 38
$ws.Range("D57").Value = 204
$ws.Range("E57").Value = 249
$ws.Range("F57").Value = 1
$ws.Range("H57").Value = 4
$ws.Range("B58").Value = 456
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 95
$ws.Range("E58").Value = 340
$ws.Range("F58").Value = 0
$ws.Range("H58").Value = 21
$ws.Range("B59").Value = 453
$ws.Range("C59").Value = 42
$ws.Range("D59").Value = 110
$ws.Range("E59").Value = 339
$ws.Range("F59").Value = 4
$ws.Range("B86").Value = 153
$ws.Range("C86").Value = 5
$ws.Range("D86").Value = 17
$ws.Range("E86").Value = 136
$ws.Range("F86").Value = 3
$ws.Range("H86").Value = 0
$ws.Range("B87").Value = 149
$ws.Range("D87").Value = 2
$ws.Range("E87").Value = 146
$ws.Range("F87").Value = 28
$ws.Range("H87").Value = 1
$ws.Range("C90").Value = 64
$ws.Range("D90").Value = 1
$ws.Range("E90").Value = 127
$ws.Range("F90").Value = 0
$ws.Range("H90").Value = 4
$ws.Range("B91").Value = 132
$ws.Range("D91").Value = 4
$ws.Range("E91").Value = 125
$ws.Range("F91").Value = 3
$ws.Range("H91").Value = 3
$ws.Range("B92").Value = 129
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 2
$ws.Range("E92").Value = 127
$ws.Range("F92").Value = 1
$ws.Range("B93").Value = 122
$ws.Range("C93").Value = 29
$ws.Range("D93").Value = 10
$ws.Range("F93").Value = 6
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 3
$ws.Range("B94").Value = 115
$ws.Range("C94").Value = 4
$ws.Range("D94").Value = 1
$ws.Range("E94").Value = 114
$ws.Range("B95").Value = 114
$ws.Range("C95").Value = 5
$ws.Range("D95").Value = 5
$ws.Range("E95").Value = 109
$ws.Range("F95").Value = 1
$ws.Range("B96").Value = 109
$ws.Range("C96").Value = 28
$ws.Range("D96").Value = 2
$ws.Range("E96").Value = 107
$ws.Range("F96").Value = 0
$ws.Range("B97").Value = 109
$ws.Range("C97").Value = 10
$ws.Range("D97").Value = 23
$ws.Range("E97").Value = 86
$ws.Range("B98").Value = 106
$ws.Range("D98").Value = 15
$ws.Range("E98").Value = 91
$ws.Range("F98").Value = 2
$ws.Range("B99").Value = 105
$ws.Range("C99").Value = 6
$ws.Range("D99").Value = 9
$ws.Range("E99").Value = 96
$ws.Range("F99").Value = 0
$ws.Range("B100").Value = 102
$ws.Range("D100").Value = 7
$ws.Range("E100").Value = 95
$ws.Range("F100").Value = 3
$ws.Range("H100").Value = 0
$ws.Range("B101").Value = 96
$ws.Range("D101").Value = 10
$ws.Range("E101").Value = 86
$ws.Range("F101").Value = 1
$ws.Range("B102").Value = 86
$ws.Range("D102").Value = 29
$ws.Range("E102").Value = 57
$ws.Range("F102").Value = 2
$ws.Range("H102").Value = 0
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 80
$ws.Range("H103").Value = 2
$ws.Range("B104").Value = 84
$ws.Range("C104").Value = 13
$ws.Range("D104").Value = 17
$ws.Range("E104").Value = 66
$ws.Range("H104").Value = 1
$ws.Range("B105").Value = 80
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 3
$ws.Range("E105").Value = 77
$ws.Range("F105").Value = 0
$ws.Range("B106").Value = 77
$ws.Range("C106").Value = 2
$ws.Range("D106").Value = 10
$ws.Range("E106").Value = 67
$ws.Range("F106").Value = 1
$ws.Range("H106").Value = 0
$ws.Range("B107").Value = 75
$ws.Range("D107").Value = 2
$ws.Range("F107").Value = 0
$ws.Range("B108").Value = 73
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 72
$ws.Range("F108").Value = 4
$ws.Range("H108").Value = 1
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 2
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 1
$ws.Range("C118").Value = 3
$ws.Range("D118").Value = 0
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 3
$ws.Range("C148").Value = 1
$ws.Range("C151").Value = 0
$ws.Range("C157").Value = 0
$ws.Range("C158").Value = 1
$ws.Range("C161").Value = 0
$ws.Range("C162").Value = 3
$ws.Range("C178").Value = 1
$ws.Range("C185").Value = 0
$ws.Range("C194").Value = 0
$ws.Range("C195").Value = 1
